# Add data for 2025-12-09
# Updates 2025 (column L) violent-crime figures across the Citywide Totals,
# By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 6220
$ws.Cells.Item(3, 12).Value = 6750
$ws.Cells.Item(4, 12).Value = 1669
$ws.Cells.Item(5, 12).Value = 399
$ws.Cells.Item(6, 12).Value = 5551
$ws.Cells.Item(7, 12).Value = 20589

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 12).Value = 182
$ws.Cells.Item(8, 12).Value = 1360
$ws.Cells.Item(15, 12).Value = 168
$ws.Cells.Item(16, 12).Value = 47
$ws.Cells.Item(17, 12).Value = 36
$ws.Cells.Item(19, 12).Value = 556
$ws.Cells.Item(20, 12).Value = 526
$ws.Cells.Item(22, 12).Value = 67
$ws.Cells.Item(23, 12).Value = 218
$ws.Cells.Item(24, 12).Value = 60
$ws.Cells.Item(25, 12).Value = 125
$ws.Cells.Item(27, 12).Value = 177
$ws.Cells.Item(29, 12).Value = 1150
$ws.Cells.Item(35, 12).Value = 26
$ws.Cells.Item(37, 12).Value = 792
$ws.Cells.Item(42, 12).Value = 657
$ws.Cells.Item(43, 12).Value = 153
$ws.Cells.Item(47, 12).Value = 146
$ws.Cells.Item(48, 12).Value = 270
$ws.Cells.Item(50, 12).Value = 100
$ws.Cells.Item(53, 12).Value = 230
$ws.Cells.Item(55, 12).Value = 219
$ws.Cells.Item(63, 12).Value = 60
$ws.Cells.Item(66, 12).Value = 58
$ws.Cells.Item(67, 12).Value = 715
$ws.Cells.Item(72, 12).Value = 83
$ws.Cells.Item(73, 12).Value = 162
$ws.Cells.Item(76, 12).Value = 319
$ws.Cells.Item(79, 12).Value = 568
$ws.Cells.Item(85, 12).Value = 1024
$ws.Cells.Item(87, 12).Value = 56
$ws.Cells.Item(90, 12).Value = 217
$ws.Cells.Item(91, 12).Value = 277
$ws.Cells.Item(94, 12).Value = 251
$ws.Cells.Item(95, 12).Value = 290
$ws.Cells.Item(96, 12).Value = 228
$ws.Cells.Item(101, 12).Value = 20589

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 12).Value = 65
$ws.Cells.Item(7, 12).Value = 228

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 309
$ws.Cells.Item(6, 12).Value = 212
$ws.Cells.Item(7, 12).Value = 1024

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 12).Value = 69
$ws.Cells.Item(3, 12).Value = 58
$ws.Cells.Item(7, 12).Value = 230

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 12).Value = 412
$ws.Cells.Item(3, 12).Value = 483
$ws.Cells.Item(7, 12).Value = 1360

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 12).Value = 104
$ws.Cells.Item(7, 12).Value = 290

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 240
$ws.Cells.Item(3, 12).Value = 281
$ws.Cells.Item(7, 12).Value = 792

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 12).Value = 280
$ws.Cells.Item(6, 12).Value = 165
$ws.Cells.Item(7, 12).Value = 715

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 12).Value = 443
$ws.Cells.Item(6, 12).Value = 281
$ws.Cells.Item(7, 12).Value = 1150

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(6, 12).Value = 110
$ws.Cells.Item(7, 12).Value = 270

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 12).Value = 201
$ws.Cells.Item(4, 12).Value = 26
$ws.Cells.Item(7, 12).Value = 556

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 12).Value = 65
$ws.Cells.Item(7, 12).Value = 319

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(3, 12).Value = 226
$ws.Cells.Item(6, 12).Value = 186
$ws.Cells.Item(7, 12).Value = 657

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 12).Value = 61
$ws.Cells.Item(7, 12).Value = 219

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(4, 12).Value = 4
$ws.Cells.Item(7, 12).Value = 60

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 12).Value = 54
$ws.Cells.Item(7, 12).Value = 218

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 12).Value = 36
$ws.Cells.Item(7, 12).Value = 277

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 12).Value = 176
$ws.Cells.Item(6, 12).Value = 152
$ws.Cells.Item(7, 12).Value = 568

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 12).Value = 164
$ws.Cells.Item(7, 12).Value = 526

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(3, 12).Value = 14
$ws.Cells.Item(7, 12).Value = 36

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(3, 12).Value = 60
$ws.Cells.Item(7, 12).Value = 251

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(5, 12).Value = 6
$ws.Cells.Item(6, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 125

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 146

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(4, 12).Value = 14
$ws.Cells.Item(7, 12).Value = 168

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(3, 12).Value = 27
$ws.Cells.Item(7, 12).Value = 100

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 58

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Cells.Item(3, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 26

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(4, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 162

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 12).Value = 59
$ws.Cells.Item(6, 12).Value = 47
$ws.Cells.Item(7, 12).Value = 182

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 12).Value = 51
$ws.Cells.Item(7, 12).Value = 177

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 12).Value = 72
$ws.Cells.Item(7, 12).Value = 217

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 12).Value = 46
$ws.Cells.Item(7, 12).Value = 153

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(3, 12).Value = 26
$ws.Cells.Item(7, 12).Value = 67

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(4, 12).Value = 13
$ws.Cells.Item(7, 12).Value = 83

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(3, 12).Value = 14
$ws.Cells.Item(7, 12).Value = 56

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(6, 12).Value = 31
$ws.Cells.Item(7, 12).Value = 47
